$d = $word.ActiveDocument

function New-PkgXml([string]$bodyXml) {
    return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + $bodyXml + '</w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# ---------------------------------------------------------------------
# 1) Title paragraph: "Phase II Writeup" -> split into two runs with a
#    spell-check proofErr wrapped around "Writeup".
# ---------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(3)
$titleBody = '<w:body><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/></w:pPr><w:r><w:t xml:space="preserve">Phase II </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Writeup</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:body>'
$titlePara.Range.InsertXML((New-PkgXml $titleBody))

# ---------------------------------------------------------------------
# 2) "Guide Through the Code:" heading -> split into three bold runs
#    with a grammar-check proofErr wrapped around "Through".
# ---------------------------------------------------------------------
$guidePara = $d.Paragraphs.Item(5)
$guideBody = '<w:body><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">Guide </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:b/></w:rPr><w:t>Through</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve"> the Code:</w:t></w:r></w:p></w:body>'
$guidePara.Range.InsertXML((New-PkgXml $guideBody))

# ---------------------------------------------------------------------
# 3) GUI paragraph: drop the _GoBack bookmark that used to sit here -
#    it will be re-created further down, inside Daniel Conroy's
#    paragraph, matching the target layout.
# ---------------------------------------------------------------------
$guiPara = $d.Paragraphs.Item(9)
$guiBody = '<w:body><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/></w:pPr><w:r><w:tab/><w:t>GUI</w:t></w:r></w:p></w:body>'
$guiPara.Range.InsertXML((New-PkgXml $guiBody))

# ---------------------------------------------------------------------
# 4) Append the new "Individual Responsibilities" block after the
#    existing "Decisions" paragraph (which is last before the sectPr).
#    We replay the Decisions paragraph verbatim (unchanged) followed
#    by all the brand new paragraphs, in a single InsertXML call.
# ---------------------------------------------------------------------
$decisionsPara = $d.Paragraphs.Item(15)

$decisionsBody = '<w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:tab/></w:r><w:r><w:t>Because the user' + [char]0x2019 + 's requests were mostly made concerning the graphics, and seeing as this was not the primary goal of the project, we thought that these issues were most applicable for the third phase rather than the second. Therefore, in this phase we focused on cohesive and encapsulated code, maximizing reuse, and providing the basic for the storage of categories, colors, and preparing the ground for basic rendering.</w:t></w:r><w:r><w:t xml:space="preserve"> Time was spent attempting to store data online, but for time purposes we decided to store in the local repository.</w:t></w:r></w:p>'

$emptyPara = '<w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/></w:pPr></w:p>'

$indivRespHeading = '<w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>Individual Responsibilities:</w:t></w:r></w:p>'

$kurtPara = '<w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:tab/></w:r><w:r><w:t>Kurt Andres</w:t></w:r></w:p>'

$danielPara = '<w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:ind w:firstLine="720"/></w:pPr><w:r><w:t>Daniel Conroy ' + [char]0x2013 + ' I implemented the Category class and its interaction with the rest of the application. I also modified timelines and events to accommodate categories</w:t></w:r><w:r><w:t>, while also refactoring timelines and events to a degree</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t>. I made refactoring changes to event labels and aided to a degree in the saving and loading of timelines.</w:t></w:r><w:r><w:t xml:space="preserve"> I participated in discussion over the inclusion and structure of categories as they relate to the rest of the data structures.</w:t></w:r></w:p>'

$leannePara = '<w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:ind w:firstLine="720"/></w:pPr><w:r><w:t>Leanne Miller</w:t></w:r></w:p>'

$andrewPara = '<w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:ind w:firstLine="720"/></w:pPr><w:r><w:t>Andrew Thompson</w:t></w:r></w:p>'

$fullBody = '<w:body>' + $decisionsBody + $emptyPara + $indivRespHeading + $kurtPara + $danielPara + $leannePara + $andrewPara + '</w:body>'

$decisionsPara.Range.InsertXML((New-PkgXml $fullBody))

# The InsertXML above targeted what used to be the very last paragraph
# in the document body, so the engine leaves one stray empty paragraph
# behind at the true end (right before the sectPr). Trim it away so
# "Andrew Thompson" is again the final paragraph, as in the target.
$n = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($n)
if ($lastPara.Range.Text -eq "") {
    $trim = $d.Range($lastPara.Range.Start - 1, $lastPara.Range.End)
    $trim.Delete()
}

Write-Output ("Final paragraph count: " + $d.Paragraphs.Count)
